# GPLIM-4825 Handle tubes w UMI and add spacer to output
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Capture the pre-existing styling we need to reuse/relocate BEFORE any
# values are overwritten, so COM can match existing style records
# instead of synthesizing brand-new (duplicate) ones.
#   B1 (old "length" header) currently carries the bold Arial/dark-gray
#       header style that B1:E1 should keep.
#   B3 (old "length" data cell, value 3) currently carries the regular
#       (non-bold) Arial/dark-gray style that C3/D3 should adopt.
# ---------------------------------------------------------------------
$headerStyleSrc = $ws.Cells.Item(1,2)
$headerStyleSrc.Copy()
foreach ($col in 2..5) {
    $ws.Cells.Item(1,$col).PasteSpecial($xlPasteFormats)
}

$dataStyleSrc = $ws.Cells.Item(3,2)
$dataStyleSrc.Copy()
foreach ($col in 3..4) {
    $ws.Cells.Item(3,$col).PasteSpecial($xlPasteFormats)
}

# A1 becomes bold Calibri (the workbook default font, just bold) instead
# of the bold Arial used before. Build that style on a scratch cell
# (so COM matches the already-existing bold+default-font style record)
# then paste it onto A1 and wipe the scratch cell again.
$scratch = $ws.Cells.Item(100,100)
$scratch.Font.Bold = $true
$scratch.Copy()
$ws.Cells.Item(1,1).PasteSpecial($xlPasteFormats)
$scratch.Clear()

$excel.CutCopyMode = $false

# ---- Header row (row 1) ----
$ws.Cells.Item(1,1).Value = "Vessel Type"
$ws.Cells.Item(1,2).Value = "Barcode"
$ws.Cells.Item(1,3).Value = "UMI Length"
$ws.Cells.Item(1,4).Value = "Spacer Length"
$ws.Cells.Item(1,5).Value = "Location"

# ---- Data rows ----
# Row 2
$ws.Cells.Item(2,1).Value = "Eppendorf96"
$ws.Cells.Item(2,2).Value = 12345
$ws.Cells.Item(2,3).Value = 6
$ws.Cells.Item(2,4).Value = 3
$ws.Cells.Item(2,5).Value = "Inline First Read"

# Row 3
$ws.Cells.Item(3,1).Value = "Eppendorf96"
$ws.Cells.Item(3,2).ClearFormats()
$ws.Cells.Item(3,2).Value = 34567
$ws.Cells.Item(3,3).Value = 3
$ws.Cells.Item(3,4).Value = 1
$ws.Cells.Item(3,5).Value = "Before Second Index Read"

# Row 4
$ws.Cells.Item(4,1).Value = "Eppendorf96"
$ws.Cells.Item(4,2).Value = 66789
$ws.Cells.Item(4,3).Value = 9
$ws.Cells.Item(4,4).Value = 1
$ws.Cells.Item(4,5).Value = "Inline Second Read"

# Row 5
$ws.Cells.Item(5,1).Value = "Eppendorf96"
$ws.Cells.Item(5,2).Value = 77891
$ws.Cells.Item(5,3).Value = 9
$ws.Cells.Item(5,4).Value = 3
$ws.Cells.Item(5,5).Value = "Inline Second Read"

# Row 6 (new)
$ws.Cells.Item(6,1).Value = "MatrixTube075"
$ws.Cells.Item(6,2).Value = 87654
$ws.Cells.Item(6,3).Value = 3
$ws.Cells.Item(6,4).Value = 2
$ws.Cells.Item(6,5).Value = "Before First Read"

# Row 7 (new)
$ws.Cells.Item(7,1).Value = "MatrixTube075"
$ws.Cells.Item(7,2).Value = 87654
$ws.Cells.Item(7,3).Value = 3
$ws.Cells.Item(7,4).Value = 2
$ws.Cells.Item(7,5).Value = "Before Second Read"

# ---- Column widths (closest achievable values to the authored widths) ----
$ws.Columns.Item(1).ColumnWidth = 11.75
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(3).ColumnWidth = 23.666666666666664
$ws.Columns.Item(4).ColumnWidth = 23.666666666666664
$ws.Columns.Item(5).ColumnWidth = 19.916666666666668

# ---- Selection matches final edit location ----
$ws.Range("E7").Select() | Out-Null
